# This edit rotates/swaps the species-observation data among a handful of
# rows (17, 18, 20, 22, 23, 25, 26, 29, 30) while leaving the shared,
# location-level fields (D, I, P, S, T, U, V, W, Y, AA, AD, AE, AG, AT, AW,
# AX, AY) untouched. Values below are the final, literal target values for
# each affected cell (derived from the target diff), so we just assign them
# directly rather than trying to "swap" via intermediate variables.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17 (becomes the old row 20 content: Tretåig hackspett) ---
$ws.Range("A17").Value = 130960378
$ws.Range("B17").Value = 57884
$ws.Range("E17").Value = 100109
$ws.Range("F17").Value = "Tretåig hackspett"
$ws.Range("G17").Value = "Picoides tridactylus"
$ws.Range("H17").Value = "(Linnaeus, 1758)"
$ws.Range("M17").Value = "äldre spår"
$ws.Range("Q17").Value = 446272
$ws.Range("R17").Value = 6759739

# --- Row 18 (becomes the old row 17 content: Garnlav) ---
$ws.Range("A18").Value = 130960789
$ws.Range("B18").Value = 79243
$ws.Range("E18").Value = 6425
$ws.Range("F18").Value = "Garnlav"
$ws.Range("G18").Value = "Alectoria sarmentosa"
$ws.Range("H18").Value = "(Ach.) Ach."
$ws.Range("M18").Value = ""
$ws.Range("Q18").Value = 446284
$ws.Range("R18").Value = 6759886

# --- Row 20 (becomes the old row 18 content: Garnlav) ---
$ws.Range("A20").Value = 130960843
$ws.Range("B20").Value = 79243
$ws.Range("E20").Value = 6425
$ws.Range("F20").Value = "Garnlav"
$ws.Range("G20").Value = "Alectoria sarmentosa"
$ws.Range("H20").Value = "(Ach.) Ach."
$ws.Range("M20").Value = ""
$ws.Range("Q20").Value = 446247
$ws.Range("R20").Value = 6759903

# --- Row 22 (becomes the old row 23 content: Vedskivlav) ---
$ws.Range("A22").Value = 130962722
$ws.Range("B22").Value = 79862
$ws.Range("E22").Value = 6453
$ws.Range("F22").Value = "Vedskivlav"
$ws.Range("G22").Value = "Hertelidea botryosa"
$ws.Range("H22").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q22").Value = 446008
$ws.Range("R22").Value = 6759948
$ws.Range("Z22").Value = "10:26"
$ws.Range("AB22").Value = "10:26"
$ws.Range("AC22").Value = ""

# --- Row 23 (becomes the old row 22 content: Garnlav) ---
$ws.Range("A23").Value = 130963976
$ws.Range("B23").Value = 79243
$ws.Range("E23").Value = 6425
$ws.Range("F23").Value = "Garnlav"
$ws.Range("G23").Value = "Alectoria sarmentosa"
$ws.Range("H23").Value = "(Ach.) Ach."
$ws.Range("Q23").Value = 445929
$ws.Range("R23").Value = 6760099
$ws.Range("Z23").Value = "14:08"
$ws.Range("AB23").Value = "14:08"
$ws.Range("AC23").Value = "Miljöbild"

# --- Row 25 (becomes the old row 26 content: Garnlav) ---
$ws.Range("A25").Value = 130962090
$ws.Range("B25").Value = 79243
$ws.Range("E25").Value = 6425
$ws.Range("F25").Value = "Garnlav"
$ws.Range("G25").Value = "Alectoria sarmentosa"
$ws.Range("H25").Value = "(Ach.) Ach."
$ws.Range("M25").Value = ""
$ws.Range("Q25").Value = 446080
$ws.Range("R25").Value = 6759960

# --- Row 26 (becomes the old row 25 content: Spillkråka) ---
$ws.Range("A26").Value = 130961746
$ws.Range("B26").Value = 57881
$ws.Range("E26").Value = 100049
$ws.Range("F26").Value = "Spillkråka"
$ws.Range("G26").Value = "Dryocopus martius"
$ws.Range("H26").Value = "(Linnaeus, 1758)"
$ws.Range("M26").Value = "färska spår"
$ws.Range("Q26").Value = 446098
$ws.Range("R26").Value = 6760061

# --- Row 29 (becomes the old row 30 content: Vedflamlav) ---
$ws.Range("A29").Value = 130962736
$ws.Range("B29").Value = 79833
$ws.Range("E29").Value = 229821
$ws.Range("F29").Value = "Vedflamlav"
$ws.Range("G29").Value = "Ramboldia elabens"
$ws.Range("H29").Value = "(Fr.) Kantvilas & Elix"
$ws.Range("M29").Value = ""
$ws.Range("Q29").Value = 446008
$ws.Range("R29").Value = 6759948
$ws.Range("Z29").Value = "10:26"
$ws.Range("AB29").Value = "10:26"

# --- Row 30 (becomes the old row 29 content: Spillkråka) ---
$ws.Range("A30").Value = 130963807
$ws.Range("B30").Value = 57881
$ws.Range("E30").Value = 100049
$ws.Range("F30").Value = "Spillkråka"
$ws.Range("G30").Value = "Dryocopus martius"
$ws.Range("H30").Value = "(Linnaeus, 1758)"
$ws.Range("M30").Value = "färska spår"
$ws.Range("Q30").Value = 445932
$ws.Range("R30").Value = 6760079
$ws.Range("Z30").Value = "14:08"
$ws.Range("AB30").Value = "14:08"
